$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) renames: shared string text updates ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Updated GDP ("C" column) predictions and a few "Colony" (AL) flags ---
$ws.Range("C2").Value = 2839.92516805933
$ws.Range("C4").Value = 1610.923908476106
$ws.Range("C5").Value = 1280.225469721551
$ws.Range("C6").Value = 5730.354774594881
$ws.Range("C7").Value = 6967.266654334572
$ws.Range("C8").Value = 2983.242707849043
$ws.Range("C9").Value = 2898.942214704482
$ws.Range("C10").Value = 1503.870423231357
$ws.Range("C11").Value = 1955.461557360978
$ws.Range("C12").Value = 11286.24301624575
$ws.Range("C13").Value = 6336.709213679884
$ws.Range("C14").Value = 4633.590358399045
$ws.Range("C15").Value = 4355.934938677345
$ws.Range("C16").Value = 5082.354756663512
$ws.Range("C17").Value = 11992.01662617741
$ws.Range("C18").Value = 2948.84548976845
$ws.Range("C20").Value = 2965.153206179127
$ws.Range("C21").Value = 1939.33862702996
$ws.Range("C22").Value = 5660.517066940175
$ws.Range("AL22").Value = 1
$ws.Range("C23").Value = 8841.561277324312
$ws.Range("C24").Value = 1577.487171555845
$ws.Range("C25").Value = 3083.80337578809
$ws.Range("C26").Value = 8390.479071096475
$ws.Range("C27").Value = 5885.254624554112
$ws.Range("C28").Value = 6947.966251196303
$ws.Range("C29").Value = 14179.19231490798
$ws.Range("C30").Value = 5360.226632400601
$ws.Range("C31").Value = 4921.848409120176
$ws.Range("C32").Value = 6711.616186806423
$ws.Range("C33").Value = 2024.117324382548
$ws.Range("C34").Value = 11627.81065059172
$ws.Range("C35").Value = 4479.398934239905
$ws.Range("C36").Value = 13455.83781255333
$ws.Range("C37").Value = 10883.31535948899
$ws.Range("C38").Value = 9477.887185090232
$ws.Range("C39").Value = 1263.452411343738
$ws.Range("C40").Value = 4022.237688257
$ws.Range("C41").Value = 1629.435089125503
$ws.Range("C42").Value = 4524.720276132375
$ws.Range("C43").Value = 12574.90356995006
$ws.Range("C44").Value = 2995.45235738661
$ws.Range("C46").Value = 3156.723844635973
$ws.Range("C47").Value = 1657.651524528445
$ws.Range("C48").Value = 2094.024217383061
$ws.Range("C49").Value = 6911.59200404802
$ws.Range("C50").Value = 5122.180090208862
$ws.Range("C51").Value = 5642.578115155247
$ws.Range("C52").Value = 11745.7759262897
$ws.Range("C53").Value = 3087.12349650562
$ws.Range("C55").Value = 3212.740625904757
$ws.Range("C56").Value = 1716.389195271215
$ws.Range("C57").Value = 2201.396847776877
$ws.Range("C58").Value = 7200.731056811853
$ws.Range("C59").Value = 5295.682695961288
$ws.Range("C60").Value = 5919.20956823756
$ws.Range("C61").Value = 11993.48398487312
$ws.Range("C62").Value = 2286.013198234259
$ws.Range("C63").Value = 11951.20944634967
$ws.Range("C64").Value = 7449.08671983612
$ws.Range("C65").Value = 1775.027517189621
$ws.Range("C66").Value = 6301.696269820412
$ws.Range("AL66").Value = 1
$ws.Range("C67").Value = 2361.056581219794
$ws.Range("C68").Value = 11431.15448084494
$ws.Range("C69").Value = 7580.275568826287
$ws.Range("C70").Value = 1836.014008604312
$ws.Range("C71").Value = 6661.86504232374
$ws.Range("AL71").Value = 1
$ws.Range("C72").Value = 7633.969039669125
$ws.Range("C73").Value = 1895.214690888655
$ws.Range("C74").Value = 7026.178156858586
$ws.Range("AL74").Value = 1

Write-Host "Edit applied successfully"
